$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRP")

$ws.Cells.Item(6, 5).Value = "2658,9"
$ws.Cells.Item(7, 5).Value = "2390,5"
$ws.Cells.Item(9, 5).Value = "1164,8"
$ws.Cells.Item(10, 5).Value = "1260,8"
$ws.Cells.Item(12, 5).Value = "2172,6"
$ws.Cells.Item(13, 5).Value = "910,7"
$ws.Cells.Item(14, 5).Value = "1890,7"
$ws.Cells.Item(15, 5).Value = "1637,7"
$ws.Cells.Item(16, 5).Value = "992,3"
$ws.Cells.Item(17, 5).Value = "1129,7"
$ws.Cells.Item(18, 5).Value = "1062,5"
$ws.Cells.Item(19, 5).Value = "1016,1"
$ws.Cells.Item(20, 5).Value = "801,3"
$ws.Cells.Item(21, 5).Value = "1357,9"
$ws.Cells.Item(23, 5).Value = "683,6"
$ws.Cells.Item(24, 5).Value = "440,9"
$ws.Cells.Item(25, 5).Value = "1008,8"
$ws.Cells.Item(26, 5).Value = "706,3"
$ws.Cells.Item(27, 5).Value = "584,4"
$ws.Cells.Item(28, 5).Value = "497,7"
$ws.Cells.Item(29, 5).Value = "301,5"
$ws.Cells.Item(31, 5).Value = "307,7"
$ws.Cells.Item(32, 5).Value = "319,1"
$ws.Cells.Item(33, 5).Value = "463,6"
$ws.Cells.Item(34, 5).Value = "452,3"
$ws.Cells.Item(35, 5).Value = "403,7"
$ws.Cells.Item(36, 5).Value = "366,6"
$ws.Cells.Item(37, 5).Value = "244,7"
$ws.Cells.Item(38, 5).Value = "277,8"
$ws.Cells.Item(39, 5).Value = "230,3"
$ws.Cells.Item(40, 5).Value = "501,8"
$ws.Cells.Item(41, 5).Value = "359,3"
$ws.Cells.Item(43, 5).Value = "804,4"
